$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.367111802101135
$ws.Range("B1").Value = 3.678926229476929
$ws.Range("C1").Value = 3.276633024215698
$ws.Range("D1").Value = 2.665881633758545
$ws.Range("E1").Value = 1.687168717384338
